$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("board")

$ws.Range("J2:Q2").Value = 1
$ws.Range("J3:Q3").Value = 1

$ws.Activate()
$ws.Range("K3").Select()
